$p = $ppt.ActivePresentation

# --- Slide 7 ("Coding Competition"): registration-password run cleanup ---
# The run containing "!" and the run containing "SpringCampus2018" carry identical
# character formatting, so PowerPoint re-merges them into a single run
# ("!SpringCampus2018") while leaving the differently-formatted trailing "!" run
# (no smtClean flag) untouched.
$s7 = $p.Slides.Item(7)
$pwdShape = $s7.Shapes.Item("Inhaltsplatzhalter 2")
$pwdRange = $pwdShape.TextFrame.TextRange
$pwdText = $pwdRange.Text
$mergeStart = $pwdText.IndexOf("!SpringCampus2018") + 1   # TextRange.Characters is 1-based
$mergeRun = $pwdRange.Characters($mergeStart, "!SpringCampus2018".Length)
$mergeRun.Text = "!SpringCampus2018"

# --- Slide 12 ("Schedule") ---
$s12 = $p.Slides.Item(12)

# Hackathon timing moved an hour earlier: 18h45-19h45 -> 17h45-18h45.
$scheduleTable = $s12.Shapes.Item("Inhaltsplatzhalter 3").Table
$tuesdayRange = $scheduleTable.Cell(3, 1).Shape.TextFrame.TextRange
$tuesdayRange.Text = $tuesdayRange.Text.Replace("18h45 – 19h45", "17h45 – 18h45")

$wednesdayRange = $scheduleTable.Cell(4, 1).Shape.TextFrame.TextRange
$wednesdayRange.Text = $wednesdayRange.Text.Replace("18h45 – 19h45", "17h45 – 18h45")

# The "Hackathons will take place..." caption textbox is no longer needed; remove it.
$s12.Shapes.Item("Textfeld 2").Delete()
